# B6-PowerPoint.pptx edit
#
# 1) Three tables (slides 14, 15, 16) get their table style switched from
#    the "Grid Table" style {F15DEC07-CA06-4109-904E-3E52A58F2A85} to the
#    "No Style, Table Grid"-ish style {0513BEAD-41C4-4F9D-9B3B-B0D4E1DDB483}.
#
# 2) The deck's design swaps from the pink/violet "Integral" palette to the
#    default blue/grey "Office" palette (the underlying theme part keeps the
#    colour/font/format scheme structure - only the 12 theme colours change).

$p = $ppt.ActivePresentation

# --- 1) Re-style the three tables -----------------------------------------
$newStyleId = "{0513BEAD-41C4-4F9D-9B3B-B0D4E1DDB483}"
foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2) Swap the theme palette from "Integral" (Red Violet) to "Office" ---
$officeColors = @(
    @(0x00, 0x00, 0x00),  # dk1
    @(0xFF, 0xFF, 0xFF),  # lt1
    @(0x44, 0x54, 0x6A),  # dk2
    @(0xE7, 0xE6, 0xE6),  # lt2
    @(0x5B, 0x9B, 0xD5),  # accent1
    @(0xED, 0x7D, 0x31),  # accent2
    @(0xA5, 0xA5, 0xA5),  # accent3
    @(0xFF, 0xC0, 0x00),  # accent4
    @(0x44, 0x72, 0xC4),  # accent5
    @(0x70, 0xAD, 0x47),  # accent6
    @(0x05, 0x63, 0xC1),  # hlink
    @(0x95, 0x4F, 0x72)   # folHlink
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $rgb = $officeColors[$i - 1]
    $r = $rgb[0]; $g = $rgb[1]; $b = $rgb[2]
    $themeColors.Item($i).RGB = $r + ($g * 256) + ($b * 65536)
}
